$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.569.28"
$ws.Range("E2").Value = "  +5.54%  "
$ws.Range("D3").Value = "3.172.14"
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("E4").Value = "  +0.10%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "402.57"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +3.77%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "110.01"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +6.37%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.549"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.38%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.04%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.618"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +4.87%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "39.13"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +4.69%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0891"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +3.58%  "
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("D13").Value = "3.669.18"
$ws.Range("E13").Value = "  +1.75%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "19.04"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "8.06"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("E16").Value = "  +6.96%  "
$ws.Range("D17").Value = "3.168.31"
$ws.Range("E17").Value = "  +2.01%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "10.53"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -3.15%  "
$ws.Range("D19").Value = "54.398.07"
$ws.Range("E19").Value = "  +4.76%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "3.29"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.87%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.0000100"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +3.50%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "12.84"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +2.45%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "72.10"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.82%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "275.67"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +3.14%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "3.28"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +4.97%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "8.08"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.09%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "27.74"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.94%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.48"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +3.23%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("E30").Value = "  -0.21%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.111"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +2.02%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "11.19"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +7.91%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0498"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +8.68%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "36.58"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("E35").Value = "  +0.94%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "50.94"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.20%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.58"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +5.56%  "
$ws.Range("E38").Value = "  -0.07%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.97"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +14.98%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "4.10"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +11.55%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "133.28"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +2.80%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.291"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.80%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.91"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "17.24"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.37%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.117"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.79%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "22.23"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.65%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.47"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.03%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.09"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").Value = "2.099.79"
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.83"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +28.39%  "
$ws.Range("B51").Value = "FlareNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/2hOSU_JYX+flarenetwork-flr"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0510"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +10.26%  "
